$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for week 17 and 18 (rows 18 and 19)
$ws.Range("B18").Value = 607
$ws.Range("B19").Value = 388

# Add new row for week 19 (row 20)
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 2
